# Rename the "4% gap" sheet to "5% gap" and update its saved selection,
# matching the author's commit "added result of running time".

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("#E#C(有上界允许4%间隙)  ")
$ws.Name = "#E#C(有上界允许5%间隙)  "

$ws.Activate()
$ws.Range("A1:J11").Select()
$excel.ActiveCell = $ws.Range("J11")
